$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (row 11): right-answer marks from 3 to 5
$ws.Range("B11").Value = 5

# Update "Total" row (row 12): total marks obtained from 75 to 125
$ws.Range("B12").Value = 125

# Update the Corr/total marks text in E12
$ws.Range("E12").Value = "125/140"
